$d = $word.ActiveDocument

# Correct "questi" -> "queste" in:
#   "...il butler possa essere interrotto solo prima di compiere una di questi."
# Locate the exact "questi." span and rewrite just the misspelled word
# ("questi" -> "queste"), leaving the surrounding sentence/runs untouched.
$rng = $d.Content
$found = $rng.Find.Execute("questi.", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if ($found) {
    $wordRange = $d.Range($rng.Start, $rng.Start + 6)
    $wordRange.Text = "queste"
}
